$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used by "full" rows (117/118): B..H, I, J, K..AD
$colsFull = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Columns used by rows 234/235/236 (no HTHG/HTAG => skip I,J)
$colsNoIJ = @("B","C","D","E","F","G","H","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $vals[$col]
    }
}

# --- Swap rows 117 and 118 (content columns B..AD); column A (id) stays fixed ---
$row117 = Get-RowValues $ws 117 $colsFull
$row118 = Get-RowValues $ws 118 $colsFull

Set-RowValues $ws 117 $colsFull $row118
Set-RowValues $ws 118 $colsFull $row117

# --- Rotate content among rows 234, 235, 236: new234 = old236, new235 = old234, new236 = old235 ---
$row234 = Get-RowValues $ws 234 $colsNoIJ
$row235 = Get-RowValues $ws 235 $colsNoIJ
$row236 = Get-RowValues $ws 236 $colsNoIJ

Set-RowValues $ws 234 $colsNoIJ $row236
Set-RowValues $ws 235 $colsNoIJ $row234
Set-RowValues $ws 236 $colsNoIJ $row235
